# "Added round 4 bids" - every remaining taken_male/taken_female cell that
# was still "false" is flipped to "true" (bids for round 4 came in, so all
# outstanding slots are now taken). This leaves no cell referencing the
# "false" shared string, so it naturally drops out of sharedStrings.xml
# when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells (column D = taken_male, column E = taken_female) whose value is
# currently the text "false" and must become the text "true".
$targets = @("D29","E30","D33","E35","D37","D38","D39","D41","E47","D49","E49","D50","E50","E52","D53","E53","D54","D55","E55","D57","E57","E58","E59","D60","E60","D61","E61","D62","E62","D63","D64","E64","D66","D67","D68","D69","E69","E70","D71","E72","D73","D74","E74","D75","E76","D77","E77","D78","D79","E80","D82","E82","D83","E83","D84","E84","E85","D86","E86","D87","E87","E88","D90","E90","D91","E91","D92","D93","D94","E94","D95","D96","E98","E99")

# Use a cell that already holds the text value "true" (shared string) as
# the copy source so the destination cells keep the same text cell type
# and shared string entry instead of Excel auto-coercing the literal
# "true"/"false" text into a boolean value.
$src = $ws.Range("E2")

foreach ($addr in $targets) {
    $src.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null   # xlPasteValues
}
$excel.CutCopyMode = 0

# Update the view state left behind by scrolling down to review the newly
# updated rows and leaving the selection on G91.
$ws.Activate()
$ws.Range("G91").Select()
$excel.ActiveWindow.ScrollRow = 86
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
